$d = $word.ActiveDocument

# --- Change 1: expand PSO practices list in "Site Lead" bullet ---
$d.Content.Find.Execute(
    "Serve as Site Lead overseeing delivery governance across all 7 PSO practices in Southeast Asia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Serve as Site Lead overseeing delivery governance across all 7 PSO practices (Data Analytics, AI/ML, Infrastructure, Security, Enterprise Architecture, Application Development, Delivery Management) in Southeast Asia",
    2) | Out-Null

# --- Change 2: expand PSO practices list in "Pioneered agentic AI" bullet ---
$d.Content.Find.Execute(
    "Pioneered agentic AI adoption across 7 PSO practices and 6 JAPAC sub-regions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pioneered agentic AI adoption across all 7 PSO practices (Data Analytics, AI/ML, Infrastructure, Security, Enterprise Architecture, Application Development, Delivery Management) and 6 JAPAC sub-regions",
    2) | Out-Null

# --- Change 3: rename "Technical Innovation & Research (Official IP):" header ---
$d.Content.Find.Execute(
    "Technical Innovation & Research (Official IP):",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Technical Innovation & Research:",
    2) | Out-Null

# --- Change 4: consolidate the 6 "Technical Innovation" bullets into 2 ---
# Locate the first of the 6 bullet paragraphs by its distinctive leading text.
$paras = $d.Paragraphs
$startIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "Research on speculative decoding*") {
        $startIdx = $i
        break
    }
}
if ($startIdx -eq -1) {
    throw "Could not locate 'Research on speculative decoding' bullet paragraph"
}

# Rewrite the first bullet's text (preserving its paragraph mark / formatting).
$p1 = $paras.Item($startIdx)
$r1 = $p1.Range
$rr1 = $d.Range($r1.Start, $r1.End - 1)
$rr1.Text = "5 Google Technical Disclosures on AI and distributed systems - UPIR (automated system synthesis, 274x speedup), FTCS (context architecture for AI agents), ARTEMIS (multi-agent debate framework), ETLC (data processing for GenAI), and LLM inference optimization (speculative decoding, custom Triton kernels)."

# Rewrite the second bullet's text (the paragraph right after, still at $startIdx + 1).
$paras = $d.Paragraphs
$p2 = $paras.Item($startIdx + 1)
$r2 = $p2.Range
$rr2 = $d.Range($r2.Start, $r2.End - 1)
$rr2.Text = "Industry-agnostic agentic AI for enterprise trust decisions. APLS self-learning + cascade routing achieving 86% cost reduction, sub-50ms latency. Won Google Cloud PSO Hackathon JAPAC, qualified for World Finals."

# Delete the remaining 4 bullet paragraphs (originally bullets 3-6 of the 6).
$paras = $d.Paragraphs
$pFirstToDelete = $paras.Item($startIdx + 2)
$pLastToDelete = $paras.Item($startIdx + 5)
$delRange = $d.Range($pFirstToDelete.Range.Start, $pLastToDelete.Range.End)
$delRange.Delete()
